# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" between "2021-Q3" and "总计", populated
#    with per-fund holding data for the 2022 Q1 snapshot.
# 2. Prepend a "2022-Q1" summary row to the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create & position the "2022-Q1" worksheet (after "2021-Q3", before
#    "总计" which naturally stays last).
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $q3Sheet)
$newSheet.Name = "2022-Q1"

# Pull header / index-column formatting from the 2021-Q3 sheet so the new
# sheet matches the look of its siblings (bold header row + bordered,
# centred A-column).
$q3Sheet.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A13").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$fundRows = @(
    @(0,  "004814", "中欧红利优享灵活配置混合A",         "22.96", "93.96", "2.72", "0.6245", 10),
    @(1,  "004815", "中欧红利优享灵活配置混合C",         "8.45",  "93.96", "2.72", "0.2298", 10),
    @(2,  "460009", "华泰柏瑞量化先行混合A",             "9.13",  "90.47", "1.06", "0.0968", 1),
    @(3,  "011448", "中银证券均衡成长混合A",             "1.41",  "90.83", "3.75", "0.0529", 5),
    @(4,  "000826", "广发中证百度百发策略100指数A",      "4.11",  "92.42", "1.04", "0.0427", 10),
    @(5,  "000827", "广发中证百度百发策略100指数E",      "4.11",  "92.42", "1.04", "0.0427", 10),
    @(6,  "005055", "华泰柏瑞量化阿尔法灵活配置混合A",   "2.53",  "89.49", "1.09", "0.0276", 6),
    @(7,  "011449", "中银证券均衡成长混合C",             "0.61",  "90.83", "3.75", "0.0229", 5),
    @(8,  "005328", "前海开源价值策略股票",               "0.37",  "92.34", "4.00", "0.0148", 5),
    @(9,  "168301", "东海祥龙灵活配置混合（LOF）",       "0.16",  "87.56", "2.41", "0.0039", 6),
    @(10, "010246", "华泰柏瑞量化先行混合C",             "0.12",  "90.47", "1.06", "0.0013", 1),
    @(11, "006532", "华泰柏瑞量化阿尔法灵活配置混合C",   "0.01",  "89.49", "1.09", "0.0001", 6)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $row[0]

    # Columns B-G hold numeric-looking text (fund code, name, size, position
    # figures) -- force text storage so e.g. "004814" keeps its leading
    # zero, then drop the temporary number format again so the cell is
    # left unstyled like the source data.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $newSheet.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 1]
        $cell.ClearFormats()
    }

    # Column H (仓位排名) is a genuine number.
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Match the bordered/centred style used by the other index cells (A3, A4).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 12
$totalSheet.Cells.Item(2, 4).Value = 1.16

# Column A is a 0-based running index, independent of the row's original
# position -- renumber the rows pushed down by the insert (were 0,1 ->
# now need to read 1,2).
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

# Restore the originally active tab ("2020-Q4") -- adding a new sheet
# shifts Excel's focus to it by default.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "2022-Q1 sheet + 总计 update applied"
